# ==========================================================================
# Edit: roll the "inputUmowa" contract template from the autumn-2024
# submission to the spring-2025 submission (dates, decision numbers,
# submission title) plus a couple of structural touch-ups.
# ==========================================================================

$d = $word.ActiveDocument

function Replace-All($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

function Replace-Once($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 1) | Out-Null
}

# 1. Appendix title + final appendix reference both gain " - wiosna 2025"
#    right before the closing curly quote.
Replace-All "Nadleśnictwie Brzeziny”" "Nadleśnictwie Brzeziny – wiosna 2025”"

Write-Host "step1 done"

# 2. Contract date: "zawarta w dniu 20 listopada 2024 roku"
#    -> "zawarta w dniu 12 lutego 2025 roku"
Replace-Once "zawarta w dniu 20 listopada 2024 roku" "zawarta w dniu 12 lutego 2025 roku"

Write-Host "step2 done"

# 3. Move the "_GoBack" bookmark from its paragraph further down (before
#    "{#osobaReprezentującaNR2}") to the second blank paragraph right
#    after the "do Regulaminu ..." heading (document paragraph #4).
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}
$targetPara = $d.Paragraphs.Item(4)
$d.Bookmarks.Add("_GoBack", $targetPara.Range) | Out-Null

Write-Host "step3 done"

# 4. Paragraf 1 ust. 1: quote + " - wiosna 2025" on the submission title,
#    updated decision numbers/dates.
Replace-Once "rozstrzygnięć Łódzkiej Submisji Drewna Cennego w Nadleśnictwie Brzeziny przeprowadzonej na podstawie:" `
             "rozstrzygnięć „Łódzkiej Submisji Drewna Cennego w Nadleśnictwie Brzeziny – wiosna 2025”, przeprowadzonej na podstawie:"

Write-Host "step4a done"

Replace-Once "Zarządzenia nr 97 Dyrektora Generalnego Lasów Państwowych z dnia 12 października 2023 roku w sprawie zasad sprzedaży drewna w PGL LP na lata 2024-2026 (znak: EM.800.3.2023)," `
             "Zarządzenia nr 120 Dyrektora Generalnego Lasów Państwowych z dnia 1 października 2024 roku w sprawie zasad sprzedaży drewna w PGL LP na lata 2025-2026 (znak: EMK.800.4.2024),"

Write-Host "step4b done"

Replace-Once "oraz Decyzji nr 46 Dyrektora Regionalnej Dyrekcji Lasów Państwowych" `
             "oraz Decyzji nr 1/2025 Dyrektora Regionalnej Dyrekcji Lasów Państwowych"

Write-Host "step4c done"

# The manual line break before "w Łodzi" is removed and "14.10.2024"
# becomes "20.01.2025"; the submission title also gains an opening curly
# quote here.
$vt = [char]11
Replace-Once ("Lasów Państwowych $vt" + "w Łodzi z dnia ") "Lasów Państwowych w Łodzi z dnia "
Write-Host "step4d done"

Replace-Once "w Łodzi z dnia 14.10.2024 roku w sprawie przeprowadzenia Łódzkiej " `
             "w Łodzi z dnia 20.01.2025 roku w sprawie przeprowadzenia „Łódzkiej "
Write-Host "step4e done"

Replace-Once "Submisji Drewna Cennego w Nadleśnictwie Brzeziny (znak: ED.804.4.4.2024)." `
             "Submisji Drewna Cennego w Nadleśnictwie Brzeziny – wiosna 2025” (znak: ED.804.4.1.2025)."
Write-Host "step4f done"

# 5. Paragraf 2 ust. 1: appendix description gets the quoted title + wiosna 2025.
Replace-Once "zakupionego na Łódzkiej Submisji Drewna Cennego w Nadleśnictwie Brzeziny) do niniejszej Umowy" `
             "zakupionego na „Łódzkiej Submisji Drewna Cennego w Nadleśnictwie Brzeziny – wiosna 2025”) do niniejszej Umowy"
Write-Host "step5 done"

# 6. All five "31.12.2024" collection / pickup / payment deadlines become
#    "28.02.2025".
Replace-All "31.12.2024" "28.02.2025"
Write-Host "step6 done"
